$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.0009005328174680471
$ws.Range("A3").Value = 0.0009005327592603862
$ws.Range("G3").Value = 77.0
$ws.Range("H3").Value = 26.0
$ws.Range("I3").Value = 2.9615390300750732
$ws.Range("A4").Value = 0.0009005327592603862
